$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text values in column D need to be forced as text
# to avoid Excel auto-converting them to numbers.
$dCells = @("D2","D3","D4","D5","D6","D7","D11","D14","D15","D16","D17","D18","D19","D21","D22","D27","D29","D30","D31","D33","D34","D36","D38","D40","D41","D42","D44","D45","D46","D47","D49","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.451.59"
$ws.Range("D3").Value = "3.800.17"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "607.58"
$ws.Range("D6").Value = "163.66"
$ws.Range("D7").Value = "3.799.61"
$ws.Range("D11").Value = "6.96"
$ws.Range("D14").Value = "35.13"
$ws.Range("D15").Value = "4.435.17"
$ws.Range("D16").Value = "3.824.92"
$ws.Range("D17").Value = "68.411.37"
$ws.Range("D18").Value = "18.08"
$ws.Range("D19").Value = "0.113"
$ws.Range("D21").Value = "462.64"
$ws.Range("D22").Value = "9.60"
$ws.Range("D27").Value = "2.11"
$ws.Range("D29").Value = "9.98"
$ws.Range("D30").Value = "3.945.35"
$ws.Range("D31").Value = "2.62"
$ws.Range("D33").Value = "7.23"
$ws.Range("D34").Value = "29.07"
$ws.Range("D36").Value = "9.05"
$ws.Range("D38").Value = "0.147"
$ws.Range("D40").Value = "0.981"
$ws.Range("D41").Value = "3.19"
$ws.Range("D42").Value = "0.999"
$ws.Range("D44").Value = "153.19"
$ws.Range("D45").Value = "0.297"
$ws.Range("D46").Value = "46.99"
$ws.Range("D47").Value = "42.95"
$ws.Range("D49").Value = "8.38"
$ws.Range("D51").Value = "26.17"

foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).ClearFormats()
}

# Other text/percentage cells (columns B, C, E) can be set directly.
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("E11").Value = "  +10.17%  "
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -6.00%  "
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +6.68%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("E40").Value = "  -1.58%  "
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -6.49%  "
